$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price cells that look numeric stay stored as text (matching source "Price" column format)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values
$ws.Range('D2').Value = '25.880.65'
$ws.Range('E2').Value = '  -2.24%  '
$ws.Range('D3').Value = '1.754.59'
$ws.Range('E3').Value = '  -4.56%  '
$ws.Range('D5').Value = '239.30'
$ws.Range('E5').Value = '  -8.06%  '
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').Value = '0.5095'
$ws.Range('E7').Value = '  -4.92%  '
$ws.Range('D8').Value = '42.34'
$ws.Range('E8').Value = '  -5.46%  '
$ws.Range('E9').Value = '  -5.54%  '
$ws.Range('D10').Value = '0.06194'
$ws.Range('E10').Value = '  -10.96%  '
$ws.Range('D11').Value = '1.748.77'
$ws.Range('E11').Value = '  -4.96%  '
$ws.Range('D12').Value = '0.06976'
$ws.Range('E12').Value = '  -2.99%  '
$ws.Range('D13').Value = '15.77'
$ws.Range('E13').Value = '  -8.51%  '
$ws.Range('D14').Value = '0.6126'
$ws.Range('E14').Value = '  -15.37%  '
$ws.Range('D15').Value = '4.535'
$ws.Range('E15').Value = '  -8.77%  '
$ws.Range('D16').Value = '77.39'
$ws.Range('E16').Value = '  -13.09%  '
$ws.Range('E17').Value = '  -0.04%  '
$ws.Range('D18').Value = '1.001'
$ws.Range('E18').Value = '  +0.01%  '
$ws.Range('D19').Value = '25.892.70'
$ws.Range('E19').Value = '  -2.23%  '
$ws.Range('D20').Value = '0.000006903'
$ws.Range('E20').Value = '  -12.25%  '
$ws.Range('D21').Value = '11.69'
$ws.Range('E21').Value = '  -14.87%  '
$ws.Range('D22').Value = '1.973.16'
$ws.Range('E22').Value = '  -5.15%  '
$ws.Range('E23').Value = '  -10.87%  '
$ws.Range('D24').Value = '5.264'
$ws.Range('E24').Value = '  -11.89%  '
$ws.Range('D25').Value = '8.251'
$ws.Range('E25').Value = '  -9.96%  '
$ws.Range('D26').Value = '138.15'
$ws.Range('E26').Value = '  -2.90%  '
$ws.Range('D27').Value = '1.493'
$ws.Range('E27').Value = '  -12.64%  '
$ws.Range('D28').Value = '15.09'
$ws.Range('E28').Value = '  -10.76%  '
$ws.Range('D29').Value = '1.820'
$ws.Range('E29').Value = '  -15.38%  '
$ws.Range('D30').Value = '103.95'
$ws.Range('E30').Value = '  -6.05%  '
$ws.Range('D31').Value = '0.08244'
$ws.Range('E31').Value = '  -7.16%  '
$ws.Range('D32').Value = '3.700'
$ws.Range('E32').Value = '  -12.74%  '
$ws.Range('D33').Value = '3.496'
$ws.Range('E33').Value = '  -13.15%  '
$ws.Range('D34').Value = '0.04554'
$ws.Range('E34').Value = '  -5.80%  '
$ws.Range('D35').Value = '0.9998'
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('D36').Value = '2.643'
$ws.Range('E36').Value = '  -8.77%  '
$ws.Range('D37').Value = '0.9934'
$ws.Range('E37').Value = '  -12.02%  '
$ws.Range('D38').Value = '0.6114'
$ws.Range('E38').Value = '  -15.49%  '
$ws.Range('D39').Value = '2.699'
$ws.Range('E39').Value = '  -12.77%  '
$ws.Range('D40').Value = '0.01561'
$ws.Range('E40').Value = '  -8.45%  '
$ws.Range('D42').Value = '103.79'
$ws.Range('E42').Value = '  -3.16%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = '1.898'
$ws.Range('E43').Value = '  -17.21%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').Value = '0.3863'
$ws.Range('E44').Value = '  -16.89%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').Value = '0.7417'
$ws.Range('E45').Value = '  -17.64%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').Value = '4.978'
$ws.Range('E46').Value = '  -15.15%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').Value = '0.05430'
$ws.Range('E47').Value = '  -5.39%  '
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').Value = '0.1115'
$ws.Range('E48').Value = '  -10.12%  '
$ws.Range('B49').Value = 'Aptos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D49').Value = '6.025'
$ws.Range('E49').Value = '  -18.54%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '7.681'
$ws.Range('E50').Value = '  -15.02%  '
$ws.Range('B51').Value = 'Elrond'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D51').Value = '30.07'
$ws.Range('E51').Value = '  -13.39%  '
